# Generate Report for Handoff
#
# Regenerates the localization-status report with a new handoff GUID
# (4710b50c-98c4-496b-a2e7-5e536c1ec974 -> 19275832-0956-4ee1-9e5f-9798ab08779b),
# a new xlf content hash (36392225dea95b570306432e4a22569d26b9958e ->
# f938f0929ce54f0f8ebaab89035e209c76f3f6bf) and refreshed handoff timestamps,
# on all three sheets (Overview, zh-cn, de-de). The hyperlink targets
# themselves (pointing at the historical git blobs) are left untouched -
# only the displayed text / cached hyperlink "display" caption changes.

$wb = $excel.ActiveWorkbook

$oldGuid = "4710b50c-98c4-496b-a2e7-5e536c1ec974"
$newGuid = "19275832-0956-4ee1-9e5f-9798ab08779b"
$oldHash = "36392225dea95b570306432e4a22569d26b9958e"
$newHash = "f938f0929ce54f0f8ebaab89035e209c76f3f6bf"

# ---------------------------------------------------------------------
# Sheet "Overview": A2 (.md hyperlink display) + D2 (plain latest-handoff
# date text, no hyperlink).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewMdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/363bc1ad7acc3b3565777dfddac0ee98bba00f58/e2e/$oldGuid.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $overviewMdAddress, "", "", "$newGuid.md")

$wsOverview.Range("D2").Value = "2016-32-18 04:32:11"

# ---------------------------------------------------------------------
# Sheet "zh-cn": A2 (.md display), B2 (".md" display, unchanged text but
# hyperlink must be recreated since Delete() clears the whole collection),
# D2 (xlf display) + E2 (plain handoff-datetime text).
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhCnMdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/363bc1ad7acc3b3565777dfddac0ee98bba00f58/e2e/$oldGuid.md"
$zhCnXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b27677a203c11e2e40eaf689608de7e548fb54c3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhCnMdAddress, "", "", "$newGuid.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), $zhCnMdAddress, "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhCnXlfAddress, "", "", "$newGuid.$newHash.zh-cn.xlf")

$wsZhCn.Range("E2").Value = "2016-03-18 04:32:09"

# ---------------------------------------------------------------------
# Sheet "de-de": A2 (.md display), B2 (".md" display, unchanged text but
# hyperlink must be recreated since Delete() clears the whole collection),
# D2 (xlf display) + E2 (plain handoff-datetime text).
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deDeMdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/363bc1ad7acc3b3565777dfddac0ee98bba00f58/e2e/$oldGuid.md"
$deDeXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c455f311194b894cbf08c68bc28e3b1edebcb194/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deDeMdAddress, "", "", "$newGuid.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), $deDeMdAddress, "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deDeXlfAddress, "", "", "$newGuid.$newHash.de-de.xlf")

$wsDeDe.Range("E2").Value = "2016-03-18 04:32:11"
